$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.712.12"
$ws.Range("E2").Value = "  -0.30%  "

$ws.Range("D3").Value = "2.732.73"

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "563.29"
$ws.Range("E5").Value = "  -1.85%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.90"
$ws.Range("E6").Value = "  +1.80%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.598"
$ws.Range("E8").Value = "  -0.46%  "

$ws.Range("E9").Value = "  +0.11%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.167"
$ws.Range("E10").Value = "  +3.97%  "

$ws.Range("E11").Value = "  +2.25%  "

$ws.Range("E12").Value = "  -0.79%  "

$ws.Range("D13").Value = "3.216.09"
$ws.Range("E13").Value = "  -0.48%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.88"
$ws.Range("E14").Value = "  +1.46%  "

$ws.Range("D15").Value = "63.539.98"
$ws.Range("E15").Value = "  -0.36%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000149"
$ws.Range("E16").Value = "  +0.01%  "

$ws.Range("D17").Value = "2.736.29"
$ws.Range("E17").Value = "  -0.55%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.52"
$ws.Range("E18").Value = "  +3.10%  "

$ws.Range("E19").Value = "  -1.05%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "354.31"
$ws.Range("E20").Value = "  -0.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.57"
$ws.Range("E21").Value = "  -2.66%  "

$ws.Range("E22").Value = "  +0.23%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.520"
$ws.Range("E23").Value = "  -3.17%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.25"
$ws.Range("E24").Value = "  -1.41%  "

$ws.Range("E25").Value = "  +0.30%  "

$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.37"
$ws.Range("E27").Value = "  -0.26%  "

$ws.Range("E28").Value = "  +1.36%  "

$ws.Range("E29").Value = "  +1.60%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.19"
$ws.Range("E30").Value = "  +3.48%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.33"
$ws.Range("E31").Value = "  +10.18%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.36"
$ws.Range("E32").Value = "  -2.30%  "

$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.02"
$ws.Range("E33").Value = "  -0.49%  "

$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.89"
$ws.Range("E34").Value = "  +1.05%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.03%  "

$ws.Range("E36").Value = "  +2.26%  "

$ws.Range("E37").Value = "  +1.20%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.972"
$ws.Range("E38").Value = "  -0.71%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "346.20"
$ws.Range("E39").Value = "  +6.73%  "

$ws.Range("E40").Value = "  +2.35%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.09"
$ws.Range("E41").Value = "  -0.94%  "

$ws.Range("E42").Value = "  -1.00%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.88"
$ws.Range("E43").Value = "  +3.07%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.03"
$ws.Range("E44").Value = "  -0.78%  "

$ws.Range("E45").Value = "  -0.55%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.627"
$ws.Range("E46").Value = "  +0.74%  "

$ws.Range("E47").Value = "  -1.44%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0999"
$ws.Range("E48").Value = "  -0.53%  "

$ws.Range("B49").Value = "FirstDigitalUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.999"
$ws.Range("E49").Value = "  +0.00%  "

$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.04"
$ws.Range("E50").Value = "  -1.78%  "

$ws.Range("E51").Value = "  -0.05%  "
